# Auto-generated script applying cached-value updates to Gilgamesh_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 173.71428
$ws.Range("I12").Value = 173.71428
$ws.Range("K12").Value = 173.71428
$ws.Range("M12").Value = -3.714280000000002
$ws.Range("H21").Value = 931.6667
$ws.Range("I21").Value = 900
$ws.Range("K21").Value = 900
$ws.Range("M21").Value = -432
$ws.Range("H23").Value = 931.6667
$ws.Range("I23").Value = 900
$ws.Range("K23").Value = 900
$ws.Range("M23").Value = -666
$ws.Range("H112").Value = 1545.3334
$ws.Range("I112").Value = 446.25
$ws.Range("J112").Value = 1765.15
$ws.Range("K112").Value = 1338.75
$ws.Range("L112").Value = 5295.450000000001
$ws.Range("M112").Value = -230.75
$ws.Range("N112").Value = -7511.450000000001
$ws.Range("H132").Value = 6355.4346
$ws.Range("I132").Value = 6355.4346
$ws.Range("K132").Value = 19066.3038
$ws.Range("M132").Value = -16536.3038
$ws.Range("H137").Value = 1319517
$ws.Range("I137").Value = 2003162.8
$ws.Range("J137").Value = 4813.5386
$ws.Range("K137").Value = 6009488.4
$ws.Range("L137").Value = 14440.6158
$ws.Range("M137").Value = -6006938.4
$ws.Range("N137").Value = -19540.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1764931.1
$ws.Range("I32").Value = 806427.0600000001
$ws.Range("K32").Value = 806427.0600000001
$ws.Range("M32").Value = -806140.0600000001
$ws.Range("H74").Value = 160821
$ws.Range("I74").Value = 266246.62
$ws.Range("J74").Value = 2682.5715
$ws.Range("K74").Value = 266246.62
$ws.Range("L74").Value = 2682.5715
$ws.Range("M74").Value = -265372.62
$ws.Range("N74").Value = -4430.5715
$ws.Range("H77").Value = 160821
$ws.Range("I77").Value = 266246.62
$ws.Range("J77").Value = 2682.5715
$ws.Range("K77").Value = 1331233.1
$ws.Range("L77").Value = 13412.8575
$ws.Range("M77").Value = -1326865.1
$ws.Range("N77").Value = -22148.8575
$ws.Range("H110").Value = 3500
$ws.Range("J110").Value = 3500
$ws.Range("L110").Value = 3500
$ws.Range("N110").Value = -7590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 699.375
$ws.Range("I22").Value = 612.75
$ws.Range("K22").Value = 612.75
$ws.Range("M22").Value = -439.75
$ws.Range("H59").Value = 111006.664
$ws.Range("J59").Value = 111006.664
$ws.Range("L59").Value = 111006.664
$ws.Range("N59").Value = -112700.664
$ws.Range("H105").Value = 20001976
$ws.Range("I105").Value = 1668614.1
$ws.Range("J105").Value = 35716284
$ws.Range("K105").Value = 1668614.1
$ws.Range("L105").Value = 35716284
$ws.Range("M105").Value = -1666867.1
$ws.Range("N105").Value = -35719778
$ws.Range("H107").Value = 4275306.5
$ws.Range("I107").Value = 5496301
$ws.Range("J107").Value = 1826.5
$ws.Range("K107").Value = 5496301
$ws.Range("L107").Value = 1826.5
$ws.Range("M107").Value = -5494381
$ws.Range("N107").Value = -5666.5
$ws.Range("H134").Value = 1185.4857
$ws.Range("I134").Value = 838.4838999999999
$ws.Range("K134").Value = 2515.4517
$ws.Range("M134").Value = 19.54830000000038

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 203.9
$ws.Range("I7").Value = 119.85714
$ws.Range("K7").Value = 119.85714
$ws.Range("M7").Value = -6.857140000000001
$ws.Range("H16").Value = 2498
$ws.Range("I16").Value = 2498
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2498
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2211
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 1896433.6
$ws.Range("I31").Value = 1756.2106
$ws.Range("J31").Value = 2662367
$ws.Range("K31").Value = 1756.2106
$ws.Range("L31").Value = 2662367
$ws.Range("M31").Value = -1461.2106
$ws.Range("N31").Value = -2662957
$ws.Range("H34").Value = 1896433.6
$ws.Range("I34").Value = 1756.2106
$ws.Range("J34").Value = 2662367
$ws.Range("K34").Value = 1756.2106
$ws.Range("L34").Value = 2662367
$ws.Range("M34").Value = -1554.2106
$ws.Range("N34").Value = -2662771
$ws.Range("H86").Value = 6941.7144
$ws.Range("I86").Value = 6606.75
$ws.Range("J86").Value = 7388.3335
$ws.Range("K86").Value = 6606.75
$ws.Range("L86").Value = 7388.3335
$ws.Range("M86").Value = -5483.75
$ws.Range("N86").Value = -9634.333500000001
$ws.Range("H89").Value = 6941.7144
$ws.Range("I89").Value = 6606.75
$ws.Range("J89").Value = 7388.3335
$ws.Range("K89").Value = 33033.75
$ws.Range("L89").Value = 36941.6675
$ws.Range("M89").Value = -27417.75
$ws.Range("N89").Value = -48173.6675
$ws.Range("H113").Value = 2498
$ws.Range("I113").Value = 2498
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2498
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -328
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1013.92
$ws.Range("I2").Value = 89.8
$ws.Range("K2").Value = 538.8
$ws.Range("M2").Value = -425.8
$ws.Range("H11").Value = 898.1667
$ws.Range("I11").Value = 898.1667
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2694.5001
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2554.5001
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 193.5
$ws.Range("I13").Value = 193.5
$ws.Range("K13").Value = 580.5
$ws.Range("M13").Value = -412.5
$ws.Range("H33").Value = 395.75
$ws.Range("I33").Value = 395
$ws.Range("J33").Value = 396.5
$ws.Range("K33").Value = 2370
$ws.Range("L33").Value = 2379
$ws.Range("M33").Value = -2087
$ws.Range("N33").Value = -2945
$ws.Range("H34").Value = 582.1667
$ws.Range("I34").Value = 373.25
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1119.75
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1035.75
$ws.Range("N34").Value = -3168
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224
$ws.Range("H39").Value = 7976.875
$ws.Range("I39").Value = 2493.3333
$ws.Range("J39").Value = 11267
$ws.Range("K39").Value = 7479.999899999999
$ws.Range("L39").Value = 33801
$ws.Range("M39").Value = -7185.999899999999
$ws.Range("N39").Value = -34389
$ws.Range("H40").Value = 246.2
$ws.Range("I40").Value = 231.14285
$ws.Range("J40").Value = 281.33334
$ws.Range("K40").Value = 924.5714
$ws.Range("L40").Value = 1125.33336
$ws.Range("M40").Value = -855.5714
$ws.Range("N40").Value = -1263.33336
$ws.Range("H55").Value = 5627.55
$ws.Range("J55").Value = 6676.7334
$ws.Range("L55").Value = 20030.2002
$ws.Range("N55").Value = -20384.2002
$ws.Range("H56").Value = 32952.6
$ws.Range("I56").Value = 32952.6
$ws.Range("K56").Value = 32952.6
$ws.Range("M56").Value = -32422.6
$ws.Range("H68").Value = 4351494
$ws.Range("I68").Value = 1052.6
$ws.Range("K68").Value = 3157.8
$ws.Range("M68").Value = -2346.8
$ws.Range("H71").Value = 4351494
$ws.Range("I71").Value = 1052.6
$ws.Range("K71").Value = 9473.4
$ws.Range("M71").Value = -5417.4
$ws.Range("H97").Value = 838083
$ws.Range("J97").Value = 8333
$ws.Range("L97").Value = 24999
$ws.Range("N97").Value = -25991
$ws.Range("H132").Value = 1609
$ws.Range("J132").Value = 1818.6666
$ws.Range("L132").Value = 16367.9994
$ws.Range("N132").Value = -21427.9994
$ws.Range("H136").Value = 8076.1816
$ws.Range("I136").Value = 1263.1428
$ws.Range("J136").Value = 19999
$ws.Range("K136").Value = 3789.4284
$ws.Range("L136").Value = 59997
$ws.Range("M136").Value = 1310.5716
$ws.Range("N136").Value = -70197
$ws.Range("H137").Value = 2646.6667
$ws.Range("J137").Value = 3382.8572
$ws.Range("L137").Value = 10148.5716
$ws.Range("N137").Value = -20348.5716
$ws.Range("H139").Value = 5364.5557
$ws.Range("I139").Value = 2305.375
$ws.Range("K139").Value = 6916.125
$ws.Range("M139").Value = -1776.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 10017.182
$ws.Range("J107").Value = 19284
$ws.Range("L107").Value = 19284
$ws.Range("N107").Value = -23124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 849.6667
$ws.Range("I61").Value = 849.6667
$ws.Range("K61").Value = 849.6667
$ws.Range("M61").Value = -647.6667
$ws.Range("H113").Value = 849.6667
$ws.Range("I113").Value = 849.6667
$ws.Range("K113").Value = 849.6667
$ws.Range("M113").Value = 1320.3333
$ws.Range("H132").Value = 4088.3157
$ws.Range("J132").Value = 4483.3335
$ws.Range("L132").Value = 13450.0005
$ws.Range("N132").Value = -18510.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("H96").Value = 1187.7778
$ws.Range("I96").Value = 1211.25
$ws.Range("K96").Value = 1211.25
$ws.Range("M96").Value = 161.75
$ws.Range("H107").Value = 943.5
$ws.Range("I107").Value = 986.4
$ws.Range("K107").Value = 2959.2
$ws.Range("M107").Value = -1039.2
$ws.Range("H122").Value = 20834712
$ws.Range("I122").Value = 1760
$ws.Range("J122").Value = 62500616
$ws.Range("K122").Value = 5280
$ws.Range("L122").Value = 187501848
$ws.Range("M122").Value = -2830
$ws.Range("N122").Value = -187506748
$ws.Range("H136").Value = 5175.8887
$ws.Range("I136").Value = 6083.4287
$ws.Range("K136").Value = 18250.2861
$ws.Range("M136").Value = -15700.2861

Write-Host "Applied updates: sets and clears complete"